$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.902.89"
$ws.Range("E2").Value = "  -3.95%  "
$ws.Range("D3").Value = "1.635.08"
$ws.Range("E3").Value = "  -6.27%  "
$ws.Range("D4").Value = "0.9970"
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").Value = "235.59"
$ws.Range("E5").Value = "  -4.44%  "
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.10%  "
$ws.Range("D7").Value = "0.4712"
$ws.Range("E7").Value = "  -6.34%  "
$ws.Range("D8").Value = "0.2545"
$ws.Range("E8").Value = "  -6.92%  "
$ws.Range("D9").Value = "0.05997"
$ws.Range("E9").Value = "  -3.00%  "
$ws.Range("D10").Value = "0.07021"
$ws.Range("E10").Value = "  -3.27%  "
$ws.Range("D11").Value = "1.635.47"
$ws.Range("E11").Value = "  -6.41%  "
$ws.Range("D12").Value = "14.73"
$ws.Range("E12").Value = "  -2.76%  "
$ws.Range("D13").Value = "0.6140"
$ws.Range("E13").Value = "  -5.96%  "
$ws.Range("D14").Value = "4.349"
$ws.Range("E14").Value = "  -6.15%  "
$ws.Range("D15").Value = "72.44"
$ws.Range("E15").Value = "  -6.70%  "
$ws.Range("D16").Value = "1.001"
$ws.Range("E16").Value = "  +0.08%  "
$ws.Range("E17").Value = "  -0.14%  "
$ws.Range("D18").Value = "24.907.57"
$ws.Range("E18").Value = "  -4.04%  "
$ws.Range("D19").Value = "0.000006545"
$ws.Range("E19").Value = "  -3.99%  "
$ws.Range("D20").Value = "11.05"
$ws.Range("E20").Value = "  -6.57%  "
$ws.Range("D21").Value = "1.846.57"
$ws.Range("E21").Value = "  -6.12%  "
$ws.Range("D22").Value = "4.369"
$ws.Range("E22").Value = "  +0.65%  "
$ws.Range("D23").Value = "8.579"
$ws.Range("E23").Value = "  -1.05%  "
$ws.Range("D24").Value = "5.256"
$ws.Range("E24").Value = "  -2.69%  "
$ws.Range("D25").Value = "133.70"
$ws.Range("E25").Value = "  -2.35%  "
$ws.Range("D26").Value = "14.76"
$ws.Range("E26").Value = "  -2.91%  "
$ws.Range("D27").Value = "1.370"
$ws.Range("E27").Value = "  -8.72%  "
$ws.Range("D28").Value = "102.30"
$ws.Range("E28").Value = "  -3.13%  "
$ws.Range("D29").Value = "1.653"
$ws.Range("E29").Value = "  -6.69%  "
$ws.Range("D30").Value = "3.743"
$ws.Range("E30").Value = "  -4.37%  "
$ws.Range("D31").Value = "0.07715"
$ws.Range("E31").Value = "  -6.39%  "
$ws.Range("D32").Value = "3.554"
$ws.Range("E32").Value = "  -2.09%  "
$ws.Range("D33").Value = "0.9991"
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("D34").Value = "0.04293"
$ws.Range("E34").Value = "  -8.19%  "
$ws.Range("D35").Value = "2.597"
$ws.Range("E35").Value = "  -2.22%  "
$ws.Range("D36").Value = "0.9168"
$ws.Range("E36").Value = "  -7.72%  "
$ws.Range("D37").Value = "0.5770"
$ws.Range("E37").Value = "  -6.67%  "
$ws.Range("D38").Value = "2.560"
$ws.Range("E38").Value = "  -6.64%  "
$ws.Range("D39").Value = "0.01547"
$ws.Range("E39").Value = "  -3.73%  "
$ws.Range("D40").Value = "0.9987"
$ws.Range("E40").Value = "  -0.13%  "
$ws.Range("D41").Value = "0.8222"
$ws.Range("E41").Value = "  +8.70%  "
$ws.Range("D42").Value = "1.793"
$ws.Range("E42").Value = "  -6.26%  "
$ws.Range("D43").Value = "96.93"
$ws.Range("E43").Value = "  -2.90%  "
$ws.Range("D44").Value = "0.3699"
$ws.Range("E44").Value = "  -4.81%  "
$ws.Range("D45").Value = "4.726"
$ws.Range("E45").Value = "  -5.59%  "
$ws.Range("D46").Value = "0.1098"
$ws.Range("E46").Value = "  -4.23%  "
$ws.Range("D47").Value = "0.05211"
$ws.Range("E47").Value = "  -0.62%  "
$ws.Range("D48").Value = "6.049"
$ws.Range("E48").Value = "  -4.01%  "
$ws.Range("D49").Value = "29.46"
$ws.Range("E49").Value = "  -3.76%  "
$ws.Range("D50").Value = "0.9996"
$ws.Range("E50").Value = "  -0.21%  "
$ws.Range("D51").Value = "0.9974"
$ws.Range("E51").Value = "  -0.73%  "
